$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Both "Manager Plugins" (row 17) and "SDN Plugins" (row 18) help-file links now
# point to the same new GUID (the new OSS-release user guide topic), replacing
# their previous distinct GUID references.
$ws.Range("C17").Value = "GUID-65309889-62B2-43BE-81CE-6A4B650AAFEE"
$ws.Range("C18").Value = "GUID-65309889-62B2-43BE-81CE-6A4B650AAFEE"

# Reflect the author's last selection/cursor position when the file was saved.
[void]$ws.Range("C13").Select()
